$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the remaining values for trial #6 (row 8) ---
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = "66.80% - Epoch 30"
$ws.Range("I8").Value = "Restarted at 16th epoch"

# --- Add a new trial #7 (row 9) ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Simple_MLP([40*(2*context_size+1), 256, 128, 71])"
$ws.Range("C9").Value = 15
$ws.Range("D9").Value = "Adam"

$ws.Range("E9").Value = 0.0002
$ws.Range("E9").NumberFormat = $ws.Range("E8").NumberFormat

$ws.Range("F9").Value = 256

$ws.Range("H8").Select()
